# Add a new worksheet "calldata" right after the existing "contacts" sheet
# and populate it with call-related lookup data, mirroring the author's edit.

$wb = $excel.ActiveWorkbook
$contacts = $wb.Worksheets.Item(1)

$calldata = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $contacts)
$calldata.Name = "calldata"

# Header row
$calldata.Range("A1").Value = "deal"
$calldata.Range("B1").Value = "task"
$calldata.Range("C1").Value = "case"
$calldata.Range("D1").Value = "note"

# Row 2
$calldata.Range("A2").Value = "a"
$calldata.Range("B2").Value = "b"
$calldata.Range("C2").Value = "c"
$calldata.Range("D2").Value = "aaa"

# Row 3
$calldata.Range("A3").Value = "e"
$calldata.Range("B3").Value = "f"
$calldata.Range("C3").Value = "g"
$calldata.Range("D3").Value = "bbb"

# Row 4 - note column filled before the rest, matching original authoring order
$calldata.Range("D4").Value = "xxx"
$calldata.Range("A4").Value = "h"
$calldata.Range("B4").Value = "i"
$calldata.Range("C4").Value = "j"

# Make the new sheet the active tab with the same selection the author left it in
$calldata.Activate()
$calldata.Range("C4").Select()
